{"js": "// Apply resume bullet-point enhancement edits (see commit message:\n// \"Enhance job descriptions across all resume types\").\n// Each entry maps an exact paragraph's current text to its replacement.\n// `occurrence` picks which match to replace when the same text appears\n// more than once in the document (1 = first match encountered, in\n// document order). Only \"Political Research and Data Analysis\" repeats\n// (under RESEARCH DIRECTOR and again under PROGRAMMER); only the first\n// (RESEARCH DIRECTOR) instance is renamed per the diff.\nconst edits = [\n  { old: \"\u2022 Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations\",\n    new: \"\u2022 Lead comprehensive research initiatives for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions\" },\n  { old: \"\u2022 Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics\",\n    new: \"\u2022 Architect enterprise-scale cloud data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics and demographic analysis\" },\n  { old: \"\u2022 Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets\",\n    new: \"\u2022 Design and implement scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets\" },\n  { old: \"\u2022 Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering\",\n    new: \"\u2022 Develop advanced analytical tools and machine learning algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering\" },\n  { old: \"\u2022 Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications\",\n    new: \"\u2022 Manage strategic client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications\" },\n  { old: \"\u2022 Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices\",\n    new: \"\u2022 Drive technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices\" },\n  { old: \"\u2022 Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES\",\n    new: \"\u2022 Conceived and developed comprehensive data framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES\" },\n  { old: \"\u2022 Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions\",\n    new: \"\u2022 Architected and built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS processing millions of records with millions of columns for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions\" },\n  { old: \"\u2022 Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI\",\n    new: \"\u2022 Led training initiatives for analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI\" },\n  { old: \"\u2022 Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company's distinguishing products\",\n    new: \"\u2022 Developed five-year strategic plans for data warehouse architecture using Scala, PySpark, and Apache Spark that became foundation of company's distinguishing products\" },\n  { old: \"\u2022 Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices\",\n    new: \"\u2022 Led cross-functional teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices\" },\n  { old: \"\u2022 Developed SimCrisis, a GeoDjango web application using Python, PostgreSQL/PostGIS, and NetLogo for multi-agent modeling and econometric simulations of crisis economies\",\n    new: \"\u2022 Architected and developed SimCrisis, a GeoDjango web application using Python, PostgreSQL/PostGIS, and NetLogo for multi-agent modeling and econometric simulations of crisis economies\" },\n  { old: \"\u2022 Liaised with officers from International Federation of Red Cross, UNICEF, and Chaos Communications Congress to improve platform using Docker and Ubuntu\",\n    new: \"\u2022 Collaborated with senior officers from International Federation of Red Cross, UNICEF, and Chaos Communications Congress to enhance platform using Docker and Ubuntu\" },\n  { old: \"\u2022 Conceived and built application using Python, Pandas, and Jupyter to predict how crisis economies respond to different humanitarian interventions\",\n    new: \"\u2022 Conceived and developed predictive application using Python, Pandas, and Jupyter to forecast how crisis economies respond to different humanitarian interventions\" },\n  { old: \"\u2022 Developed RACSO, a web application for pollsters to fully administer research including questionnaire creation, versioning, and reporting\",\n    new: \"\u2022 Architected and developed RACSO, a comprehensive web application for pollsters to fully administer research including questionnaire creation, versioning, and reporting\" },\n  { old: \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors before selecting implementation partner\",\n    new: \"\u2022 Led RFP process and analyzed bids from 1,200 vendors before selecting optimal implementation partner\" },\n  { old: \"\u2022 Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research\",\n    new: \"\u2022 Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research affecting millions of dollars in campaign spending decisions\" },\n  { old: \"\u2022 Designed survey deployment system facilitating thousands of simultaneous phone surveys\",\n    new: \"\u2022 Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly $1 million annually in polling costs\" },\n  { old: \"\u2022 Maintained and extended entire geospatial analysis and reporting tools for Java-based CRM system\",\n    new: \"\u2022 Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system\" },\n  { old: \"\u2022 Built geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill\",\n    new: \"\u2022 Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill\" },\n  { old: \"\u2022 Assisted in search for full-time CTO while performing all programmatic technology roles for multi-million dollar organization\",\n    new: \"\u2022 Led technology operations for multi-million dollar organization while assisting in search for full-time CTO\" },\n  { old: \"\u2022 Made all technology decisions and practices for massive multinational non-governmental organization\",\n    new: \"\u2022 Directed all technology decisions and practices for massive multinational non-governmental organization\" },\n  { old: \"\u2022 Wrote comprehensive frameworks for internal and external technology audits\",\n    new: \"\u2022 Developed comprehensive frameworks for internal and external technology audits\" },\n  { old: \"\u2022 Trained beneficiaries on spatial and Census data analysis for public health research\",\n    new: \"\u2022 Led training initiatives for beneficiaries on spatial and Census data analysis for public health research\" },\n  { old: \"\u2022 Trained NGO staff in web development using Drupal, PHP, and MySQL\",\n    new: \"\u2022 Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL\" },\n  { old: \"Political Research and Data Analysis\",\n    new: \"Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns\",\n    occurrence: 1 },\n  { old: \"\u2022 Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections\",\n    new: \"\u2022 Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections affecting millions of dollars in campaign spending decisions\" },\n  { old: \"Political Field Operations and Data Management\",\n    new: \"Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns\" },\n  { old: \"\u2022 Administered all quantitative and qualitative research operations ensuring reporting accuracy\",\n    new: \"\u2022 Administered all quantitative and qualitative research operations for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in spending decisions\" },\n  { old: \"\u2022 Managed comprehensive survey fielding for multi-million dollar research firm\",\n    new: \"\u2022 Managed team of 6 research analysts and field staff for comprehensive survey fielding at multi-million dollar research firm\" },\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Track how many times each \"old\" string has already been seen, so the\n// `occurrence`-limited edits (the duplicated heading) only fire on the\n// intended match.\nconst seenCounts = new Map();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  for (const edit of edits) {\n    if (text !== edit.old) continue;\n    const seenSoFar = (seenCounts.get(edit.old) || 0) + 1;\n    seenCounts.set(edit.old, seenSoFar);\n    const targetOccurrence = edit.occurrence || 1;\n    if (seenSoFar !== targetOccurrence) continue;\n    para.insertText(edit.new, \"Replace\");\n    break;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply resume bullet-point enhancement edits (see commit message:\n# \"Enhance job descriptions across all resume types\").\n# Each entry maps an exact paragraph's current text to its replacement.\n# Occurrence picks which match to replace when the same text appears\n# more than once in the document (1 = first match encountered, in\n# document order). Only \"Political Research and Data Analysis\" repeats\n# (under RESEARCH DIRECTOR and again under PROGRAMMER); only the first\n# (RESEARCH DIRECTOR) instance is renamed per the diff, so it gets\n# occurrence = 1 while the PROGRAMMER instance is left untouched.\n\n$d = $word.ActiveDocument\n\n$edits = @(\n    @{ old = \"\u2022 Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations\"; new = \"\u2022 Lead comprehensive research initiatives for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions\"; occurrence = 1 },\n    @{ old = \"\u2022 Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics\"; new = \"\u2022 Architect enterprise-scale cloud data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics and demographic analysis\"; occurrence = 1 },\n    @{ old = \"\u2022 Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets\"; new = \"\u2022 Design and implement scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets\"; occurrence = 1 },\n    @{ old = \"\u2022 Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering\"; new = \"\u2022 Develop advanced analytical tools and machine learning algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering\"; occurrence = 1 },\n    @{ old = \"\u2022 Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications\"; new = \"\u2022 Manage strategic client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications\"; occurrence = 1 },\n    @{ old = \"\u2022 Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices\"; new = \"\u2022 Drive technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices\"; occurrence = 1 },\n    @{ old = \"\u2022 Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES\"; new = \"\u2022 Conceived and developed comprehensive data framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES\"; occurrence = 1 },\n    @{ old = \"\u2022 Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions\"; new = \"\u2022 Architected and built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS processing millions of records with millions of columns for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions\"; occurrence = 1 },\n    @{ old = \"\u2022 Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI\"; new = \"\u2022 Led training initiatives for analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI\"; occurrence = 1 },\n    @{ old = \"\u2022 Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company's distinguishing products\"; new = \"\u2022 Developed five-year strategic plans for data warehouse architecture using Scala, PySpark, and Apache Spark that became foundation of company's distinguishing products\"; occurrence = 1 },\n    @{ old = \"\u2022 Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices\"; new = \"\u2022 Led cross-functional teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices\"; occurrence = 1 },\n    @{ old = \"\u2022 Developed SimCrisis, a GeoDjango web application using Python, PostgreSQL/PostGIS, and NetLogo for multi-agent modeling and econometric simulations of crisis economies\"; new = \"\u2022 Architected and developed SimCrisis, a GeoDjango web application using Python, PostgreSQL/PostGIS, and NetLogo for multi-agent modeling and econometric simulations of crisis economies\"; occurrence = 1 },\n    @{ old = \"\u2022 Liaised with officers from International Federation of Red Cross, UNICEF, and Chaos Communications Congress to improve platform using Docker and Ubuntu\"; new = \"\u2022 Collaborated with senior officers from International Federation of Red Cross, UNICEF, and Chaos Communications Congress to enhance platform using Docker and Ubuntu\"; occurrence = 1 },\n    @{ old = \"\u2022 Conceived and built application using Python, Pandas, and Jupyter to predict how crisis economies respond to different humanitarian interventions\"; new = \"\u2022 Conceived and developed predictive application using Python, Pandas, and Jupyter to forecast how crisis economies respond to different humanitarian interventions\"; occurrence = 1 },\n    @{ old = \"\u2022 Developed RACSO, a web application for pollsters to fully administer research including questionnaire creation, versioning, and reporting\"; new = \"\u2022 Architected and developed RACSO, a comprehensive web application for pollsters to fully administer research including questionnaire creation, versioning, and reporting\"; occurrence = 1 },\n    @{ old = \"\u2022 Wrote RFP and analyzed bids from 1,200 vendors before selecting implementation partner\"; new = \"\u2022 Led RFP process and analyzed bids from 1,200 vendors before selecting optimal implementation partner\"; occurrence = 1 },\n    @{ old = \"\u2022 Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research\"; new = \"\u2022 Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research affecting millions of dollars in campaign spending decisions\"; occurrence = 1 },\n    @{ old = \"\u2022 Designed survey deployment system facilitating thousands of simultaneous phone surveys\"; new = \"\u2022 Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly `$1 million annually in polling costs\"; occurrence = 1 },\n    @{ old = \"\u2022 Maintained and extended entire geospatial analysis and reporting tools for Java-based CRM system\"; new = \"\u2022 Maintained and extended comprehensive geospatial analysis and reporting tools for Java-based CRM system\"; occurrence = 1 },\n    @{ old = \"\u2022 Built geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill\"; new = \"\u2022 Built advanced geospatial analysis capabilities using Java, JavaScript, MySQL, and TileMill\"; occurrence = 1 },\n    @{ old = \"\u2022 Assisted in search for full-time CTO while performing all programmatic technology roles for multi-million dollar organization\"; new = \"\u2022 Led technology operations for multi-million dollar organization while assisting in search for full-time CTO\"; occurrence = 1 },\n    @{ old = \"\u2022 Made all technology decisions and practices for massive multinational non-governmental organization\"; new = \"\u2022 Directed all technology decisions and practices for massive multinational non-governmental organization\"; occurrence = 1 },\n    @{ old = \"\u2022 Wrote comprehensive frameworks for internal and external technology audits\"; new = \"\u2022 Developed comprehensive frameworks for internal and external technology audits\"; occurrence = 1 },\n    @{ old = \"\u2022 Trained beneficiaries on spatial and Census data analysis for public health research\"; new = \"\u2022 Led training initiatives for beneficiaries on spatial and Census data analysis for public health research\"; occurrence = 1 },\n    @{ old = \"\u2022 Trained NGO staff in web development using Drupal, PHP, and MySQL\"; new = \"\u2022 Conducted training programs for NGO staff in web development using Drupal, PHP, and MySQL\"; occurrence = 1 },\n    @{ old = \"Political Research and Data Analysis\"; new = \"Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns\"; occurrence = 1 },\n    @{ old = \"\u2022 Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections\"; new = \"\u2022 Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections affecting millions of dollars in campaign spending decisions\"; occurrence = 1 },\n    @{ old = \"Political Field Operations and Data Management\"; new = \"Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns\"; occurrence = 1 },\n    @{ old = \"\u2022 Administered all quantitative and qualitative research operations ensuring reporting accuracy\"; new = \"\u2022 Administered all quantitative and qualitative research operations for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in spending decisions\"; occurrence = 1 },\n    @{ old = \"\u2022 Managed comprehensive survey fielding for multi-million dollar research firm\"; new = \"\u2022 Managed team of 6 research analysts and field staff for comprehensive survey fielding at multi-million dollar research firm\"; occurrence = 1 }\n)\n\n# Track how many times each \"old\" string has already been seen, so the\n# occurrence-limited edits (the duplicated heading) only fire on the\n# intended match.\n$seenCounts = @{}\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd()\n    foreach ($edit in $edits) {\n        if ($text -cne $edit.old) { continue }\n        if (-not $seenCounts.ContainsKey($edit.old)) { $seenCounts[$edit.old] = 0 }\n        $seenCounts[$edit.old] = $seenCounts[$edit.old] + 1\n        if ($seenCounts[$edit.old] -ne $edit.occurrence) { continue }\n        $p.Range.Text = $edit.new\n        break\n    }\n}\n"}
